$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (the second paragraph in the document).
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Giant's Gold Free - Review and Ratings 2021"
#    right before the final "Prompt: ..." paragraph.
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs($count - 1)
$insertPoint = $secondLast.Range.Duplicate()
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("Play Giant's Gold Free - Review and Ratings 2021" + [char]13)

$newPara = $d.Paragraphs($count)
$newParaText = $newPara.Range.Duplicate()
$newParaText.MoveEnd(1, -1)
$newParaText.Font.Bold = $true

# 3. Replace the old "Prompt: ..." paragraph text with the meta-description text,
#    keeping the italic formatting.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $lastPara.Range.Duplicate()
$lastRange.MoveEnd(1, -1)
$lastRange.Text = "Play Giant's Gold free and read our review with info on gameplay mechanics, graphics, special symbols and target audience in 2021."
